# Daily attendance processing - 2026-01-23 07:50:21
# Normalize the "Recorded By" column (G): the automated "System" /
# "system" recorder tag(s) should be listed AFTER any real user/email
# entries in the comma-separated list, instead of before them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null) {
        $s = [string]$val

        if ($s.Contains(",")) {
            $parts = $s.Split(",")
            $trimmed = @()
            foreach ($p in $parts) { $trimmed += $p.Trim() }

            $others = @()
            $systemEntries = @()
            foreach ($p in $trimmed) {
                if ($p.ToLower() -eq "system") {
                    $systemEntries += $p
                } else {
                    $others += $p
                }
            }

            $newParts = $others + $systemEntries
            $newVal = [string]::Join(", ", $newParts)

            if ($newVal -ne $s) {
                $cell.Value = $newVal
            }
        }
    }
}
